$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - header/summary info row
$ws.Range("A2").Value = "10ª PmJ Mossoró"
$ws.Range("B2").Value = "Adriana Lira Da Luz Mello"
$ws.Range("C2").Value = 1
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "18-10-2024"
$ws.Range("E2").Value = "10ª PmJ Mossoró"

# Row 5
$ws.Range("B5").Value = "Último Relatório de Correição"
$ws.Range("C5").Value = "fls. 12-37"

# Row 6
$ws.Range("C6").Value = "fls. 56-57"

# Row 7
$ws.Range("C7").Value = "fls. 9"

# Row 9
$ws.Range("C9").Value = "fls. 9"

# Row 10
$ws.Range("C10").Value = "fls. 12-37"

# Row 18
$ws.Range("C18").Value = "fls. 7"

# Row 19
$ws.Range("A19").Value = "Cursos Oficiais Diversos dos de Formação Continuada"
$ws.Range("C19").Value = "fls. 7"

# Row 20
$ws.Range("A20").Value = "Cursos Reconhecidos de Aperfeiçoamento"
